$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.662.21'
$ws.Range("E2").Value = '  -3.12%  '

# Row 3
$ws.Range("D3").Value = '1.741.66'
$ws.Range("E3").Value = '  -5.40%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  +0.32%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.86'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -9.02%  '

# Row 6
$ws.Range("E6").Value = '  +0.19%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5026'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -5.99%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.73'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -6.67%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2613'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -13.64%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06125'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -11.11%  '

# Row 11
$ws.Range("D11").Value = '1.749.09'
$ws.Range("E11").Value = '  -5.60%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06947'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  -7.73%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.10'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  -16.13%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.470'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -10.37%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5904'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  -19.71%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.43'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -14.75%  '

# Row 17
$ws.Range("E17").Value = '  +0.26%  '

# Row 18
$ws.Range("E18").Value = '  +0.13%  '

# Row 19
$ws.Range("D19").Value = '25.730.83'
$ws.Range("E19").Value = '  -2.96%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.57'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -17.31%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006738'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  -15.04%  '

# Row 22
$ws.Range("D22").Value = '1.970.57'
$ws.Range("E22").Value = '  -5.61%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.036'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -12.30%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.072'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -13.25%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.080'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -15.09%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '138.17'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -3.49%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.534'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -9.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.805'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -18.73%  '

# Row 29
$ws.Range("E29").Value = '  -12.08%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.95'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -6.91%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.759'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -11.97%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08101'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -8.02%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.437'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -15.37%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04480'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -6.68%  '

# Row 35
$ws.Range("E35").Value = '  +0.09%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.636'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -9.93%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9670'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -14.71%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6013'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -17.36%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.644'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -14.75%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01542'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -10.22%  '

# Row 41
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.002'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +0.07%  '

# Row 42
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '104.02'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -3.50%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.901'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -17.20%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.119'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -12.85%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3770'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -20.59%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7249'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -19.88%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05330'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -8.07%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1102'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -10.82%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.88'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -14.50%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.849'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -21.47%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.16'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -13.65%  '
